# professores.xlsx - update CH MAX values for two teachers and refresh
# the selection / column B width like Excel does after a manual edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("professores")

# Row 12 (saulomarques.morais@gmail.com): CH MAX 4 -> 6
$ws.Range("B12").Value = 6

# Row 13 (leonardo.lima@palmares.ifpe.edu.br): CH MAX 4 -> 10
$ws.Range("B13").Value = 10

# Column B ("CH MAX") now needs to show its content, so auto-fit it
# (matches Excel's best-fit width of ~7.66 characters for this column).
$ws.Columns("B").ColumnWidth = 6.83

# Move / leave the active selection on C8 (as last saved in the workbook).
$ws.Range("C8").Select()
